# Remove the @BeforeSuite -driven "N" (negative) test row by flipping it to
# "Y" and fill in the previously-empty "browser" column for the TestData
# driver sheet (rows 7-9), then leave the selection on the newly edited range.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestData")

# Row 9 was the only row using "Runmode" = "N"; flip it to "Y" so every
# scenario now runs (removes the now-unused "N" shared string).
$ws.Range("A9").Value = "Y"

# Populate the "browser" column (D) for the second data block (rows 7-9),
# which previously had no D values.
$ws.Range("D7").Value = "browser"
$ws.Range("D8").Value = "chrome"
$ws.Range("D9").Value = "firefox"

# Match the saved selection/active cell on the TestData sheet.
$ws.Range("D7:D9").Select()
